$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1875.6923
$ws.Range("I2").Value = 1577.2
$ws.Range("K2").Value = 1577.2
$ws.Range("M2").Value = -1464.2
$ws.Range("H4").Value = 949.375
$ws.Range("I4").Value = 198.88889
$ws.Range("K4").Value = 198.88889
$ws.Range("M4").Value = -84.88889
$ws.Range("H18").Value = 1567.2858
$ws.Range("I18").Value = 1645.5
$ws.Range("J18").Value = 1536
$ws.Range("K18").Value = 1645.5
$ws.Range("L18").Value = 1536
$ws.Range("M18").Value = -1361.5
$ws.Range("N18").Value = -2104
$ws.Range("H19").Value = 2565123
$ws.Range("I19").Value = 829.61536
$ws.Range("J19").Value = 3847269.5
$ws.Range("K19").Value = 829.61536
$ws.Range("L19").Value = 3847269.5
$ws.Range("M19").Value = -654.61536
$ws.Range("N19").Value = -3847619.5
$ws.Range("H28").Value = 4522.9
$ws.Range("I28").Value = 796.9231
$ws.Range("K28").Value = 796.9231
$ws.Range("M28").Value = -311.9231
$ws.Range("H64").Value = 5306
$ws.Range("I64").Value = 4613.6665
$ws.Range("J64").Value = 6344.5
$ws.Range("K64").Value = 4613.6665
$ws.Range("L64").Value = 6344.5
$ws.Range("M64").Value = -4365.6665
$ws.Range("N64").Value = -6840.5
$ws.Range("H67").Value = 5306
$ws.Range("I67").Value = 4613.6665
$ws.Range("J67").Value = 6344.5
$ws.Range("K67").Value = 4613.6665
$ws.Range("L67").Value = 6344.5
$ws.Range("M67").Value = -3755.6665
$ws.Range("N67").Value = -8060.5
$ws.Range("H80").Value = 2586675
$ws.Range("I80").Value = 1635070.9
$ws.Range("J80").Value = 3402335.5
$ws.Range("K80").Value = 4905212.699999999
$ws.Range("L80").Value = 10207006.5
$ws.Range("M80").Value = -4904214.699999999
$ws.Range("N80").Value = -10209002.5
$ws.Range("H83").Value = 2586675
$ws.Range("I83").Value = 1635070.9
$ws.Range("J83").Value = 3402335.5
$ws.Range("K83").Value = 14715638.1
$ws.Range("L83").Value = 30621019.5
$ws.Range("M83").Value = -14710646.1
$ws.Range("N83").Value = -30631003.5
$ws.Range("H86").Value = 7853.1816
$ws.Range("I86").Value = 4326.4
$ws.Range("J86").Value = 10792.167
$ws.Range("K86").Value = 4326.4
$ws.Range("L86").Value = 10792.167
$ws.Range("M86").Value = -3203.4
$ws.Range("N86").Value = -13038.167
$ws.Range("H89").Value = 7853.1816
$ws.Range("I89").Value = 4326.4
$ws.Range("J89").Value = 10792.167
$ws.Range("K89").Value = 21632
$ws.Range("L89").Value = 53960.835
$ws.Range("M89").Value = -16016
$ws.Range("N89").Value = -65192.835
$ws.Range("H92").Value = 3110.3044
$ws.Range("J92").Value = 5216.909
$ws.Range("L92").Value = 5216.909
$ws.Range("N92").Value = -7712.909
$ws.Range("H96").Value = 1118981.9
$ws.Range("I96").Value = 3285.3333
$ws.Range("J96").Value = 2075293.1
$ws.Range("K96").Value = 9855.999899999999
$ws.Range("L96").Value = 6225879.300000001
$ws.Range("M96").Value = -8482.999899999999
$ws.Range("N96").Value = -6228625.300000001
$ws.Range("H100").Value = 4007.5
$ws.Range("I100").Value = 2231.5334
$ws.Range("J100").Value = 9335.4
$ws.Range("K100").Value = 2231.5334
$ws.Range("L100").Value = 9335.4
$ws.Range("M100").Value = -1690.5334
$ws.Range("N100").Value = -10417.4
$ws.Range("H103").Value = 125000750
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 166667330
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 500001990
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -500003162
$ws.Range("H106").Value = 10041.25
$ws.Range("I106").Value = 10041.25
$ws.Range("K106").Value = 10041.25
$ws.Range("M106").Value = -9410.25
$ws.Range("H107").Value = 782.55554
$ws.Range("I107").Value = 283.14285
$ws.Range("J107").Value = 2530.5
$ws.Range("K107").Value = 283.14285
$ws.Range("L107").Value = 2530.5
$ws.Range("M107").Value = 1636.85715
$ws.Range("N107").Value = -6370.5
$ws.Range("H111").Value = 2624.875
$ws.Range("I111").Value = 2624.875
$ws.Range("K111").Value = 7874.625
$ws.Range("M111").Value = -4807.625
$ws.Range("H132").Value = 3356.24
$ws.Range("I132").Value = 2390.9546
$ws.Range("K132").Value = 7172.8638
$ws.Range("M132").Value = -4642.8638
$ws.Range("H137").Value = 566660.9399999999
$ws.Range("I137").Value = 1394.7858
$ws.Range("J137").Value = 1445963.9
$ws.Range("K137").Value = 4184.357400000001
$ws.Range("L137").Value = 4337891.699999999
$ws.Range("M137").Value = -1634.357400000001
$ws.Range("N137").Value = -4342991.699999999
$ws.Range("H138").Value = 3055.9604
$ws.Range("I138").Value = 1830.3243
$ws.Range("J138").Value = 4218.7437
$ws.Range("K138").Value = 5490.9729
$ws.Range("L138").Value = 12656.2311
$ws.Range("M138").Value = -350.9728999999998
$ws.Range("N138").Value = -22936.2311

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1624.8
$ws.Range("I2").Value = 1500.3125
$ws.Range("K2").Value = 1500.3125
$ws.Range("M2").Value = -1387.3125
$ws.Range("H4").Value = 2931.25
$ws.Range("I4").Value = 3321.4285
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 3321.4285
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -3205.4285
$ws.Range("N4").Value = -432
$ws.Range("H13").Value = 700170
$ws.Range("I13").Value = 875075
$ws.Range("K13").Value = 875075
$ws.Range("M13").Value = -874931
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9827
$ws.Range("H19").Value = 999
$ws.Range("J19").Value = 999
$ws.Range("L19").Value = 999
$ws.Range("N19").Value = -1457
$ws.Range("H25").Value = 2993.3333
$ws.Range("J25").Value = 990
$ws.Range("L25").Value = 990
$ws.Range("N25").Value = -1794
$ws.Range("H32").Value = 3847.4922
$ws.Range("I32").Value = 3920.3386
$ws.Range("K32").Value = 3920.3386
$ws.Range("M32").Value = -3633.3386
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("K38").Value = 10000
$ws.Range("M38").Value = -9533
$ws.Range("H42").Value = 89999
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H53").Value = 17344.25
$ws.Range("J53").Value = 19826
$ws.Range("L53").Value = 19826
$ws.Range("N53").Value = -21190
$ws.Range("H61").Value = 3446283.8
$ws.Range("I61").Value = 3777006
$ws.Range("J61").Value = 1255249.9
$ws.Range("K61").Value = 3777006
$ws.Range("L61").Value = 1255249.9
$ws.Range("M61").Value = -3776794
$ws.Range("N61").Value = -1255673.9
$ws.Range("H74").Value = 783909.75
$ws.Range("I74").Value = 835503.8
$ws.Range("K74").Value = 835503.8
$ws.Range("M74").Value = -834629.8
$ws.Range("H77").Value = 783909.75
$ws.Range("I77").Value = 835503.8
$ws.Range("K77").Value = 4177519
$ws.Range("M77").Value = -4173151
$ws.Range("H97").Value = 1936.579
$ws.Range("J97").Value = 2219.6
$ws.Range("L97").Value = 2219.6
$ws.Range("N97").Value = -3211.6
$ws.Range("H109").Value = 21299.5
$ws.Range("J109").Value = 21299.5
$ws.Range("L109").Value = 21299.5
$ws.Range("N109").Value = -24073.5
$ws.Range("H116").Value = 1624.8
$ws.Range("I116").Value = 1500.3125
$ws.Range("K116").Value = 1500.3125
$ws.Range("M116").Value = 793.6875
$ws.Range("H122").Value = 2891.72
$ws.Range("I122").Value = 2419.1428
$ws.Range("K122").Value = 7257.428400000001
$ws.Range("M122").Value = -4807.428400000001
$ws.Range("H132").Value = 1153263.5
$ws.Range("I132").Value = 3709.192
$ws.Range("J132").Value = 7147368
$ws.Range("K132").Value = 11127.576
$ws.Range("L132").Value = 21442104
$ws.Range("M132").Value = -8597.576000000001
$ws.Range("N132").Value = -21447164
$ws.Range("H136").Value = 3446283.8
$ws.Range("I136").Value = 3777006
$ws.Range("J136").Value = 1255249.9
$ws.Range("K136").Value = 11331018
$ws.Range("L136").Value = 3765749.7
$ws.Range("M136").Value = -11328468
$ws.Range("N136").Value = -3770849.7

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1624.8
$ws.Range("I3").Value = 1500.3125
$ws.Range("K3").Value = 1500.3125
$ws.Range("M3").Value = -1386.3125
$ws.Range("H86").Value = 3724.875
$ws.Range("I86").Value = 2920.5
$ws.Range("K86").Value = 2920.5
$ws.Range("M86").Value = -1797.5
$ws.Range("H89").Value = 3724.875
$ws.Range("I89").Value = 2920.5
$ws.Range("K89").Value = 14602.5
$ws.Range("M89").Value = -8986.5
$ws.Range("H105").Value = 998442.4
$ws.Range("I105").Value = 1635884.4
$ws.Range("K105").Value = 1635884.4
$ws.Range("M105").Value = -1634137.4
$ws.Range("H107").Value = 5743.857
$ws.Range("I107").Value = 9883
$ws.Range("J107").Value = 3444.3333
$ws.Range("K107").Value = 9883
$ws.Range("L107").Value = 3444.3333
$ws.Range("M107").Value = -7963
$ws.Range("N107").Value = -7284.3333
$ws.Range("H134").Value = 2859381.2
$ws.Range("I134").Value = 1931.7916
$ws.Range("J134").Value = 9093816
$ws.Range("K134").Value = 5795.3748
$ws.Range("L134").Value = 27281448
$ws.Range("M134").Value = -3260.3748
$ws.Range("N134").Value = -27286518
$ws.Range("H140").Value = 99999.5
$ws.Range("J140").Value = 99999.5
$ws.Range("L140").Value = 99999.5
$ws.Range("N140").Value = -110359.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34834210
$ws.Range("I31").Value = 47621556
$ws.Range("J31").Value = 1267424.5
$ws.Range("K31").Value = 47621556
$ws.Range("L31").Value = 1267424.5
$ws.Range("M31").Value = -47621261
$ws.Range("N31").Value = -1268014.5
$ws.Range("H34").Value = 34834210
$ws.Range("I34").Value = 47621556
$ws.Range("J34").Value = 1267424.5
$ws.Range("K34").Value = 47621556
$ws.Range("L34").Value = 1267424.5
$ws.Range("M34").Value = -47621354
$ws.Range("N34").Value = -1267828.5
$ws.Range("H58").Value = 2198.1785
$ws.Range("I58").Value = 1968.25
$ws.Range("J58").Value = 3577.75
$ws.Range("K58").Value = 1968.25
$ws.Range("L58").Value = 3577.75
$ws.Range("M58").Value = -1765.25
$ws.Range("N58").Value = -3983.75
$ws.Range("H62").Value = 4691.75
$ws.Range("I62").Value = 4959.2
$ws.Range("J62").Value = 3889.4
$ws.Range("K62").Value = 4959.2
$ws.Range("L62").Value = 3889.4
$ws.Range("M62").Value = -4335.2
$ws.Range("N62").Value = -5137.4
$ws.Range("H65").Value = 4691.75
$ws.Range("I65").Value = 4959.2
$ws.Range("J65").Value = 3889.4
$ws.Range("K65").Value = 24796
$ws.Range("L65").Value = 19447
$ws.Range("M65").Value = -21676
$ws.Range("N65").Value = -25687
$ws.Range("H99").Value = 16264.5
$ws.Range("J99").Value = 26496.857
$ws.Range("L99").Value = 26496.857
$ws.Range("N99").Value = -29492.857
$ws.Range("H122").Value = 2423.923
$ws.Range("I122").Value = 2154.4211
$ws.Range("K122").Value = 6463.263300000001
$ws.Range("M122").Value = -4013.263300000001
$ws.Range("H126").Value = 16264.5
$ws.Range("J126").Value = 26496.857
$ws.Range("L126").Value = 79490.571
$ws.Range("N126").Value = -84430.571
$ws.Range("H132").Value = 1814.2941
$ws.Range("I132").Value = 1801.8125
$ws.Range("K132").Value = 5405.4375
$ws.Range("M132").Value = -2875.4375
$ws.Range("H134").Value = 2286.3333
$ws.Range("I134").Value = 2101.6155
$ws.Range("J134").Value = 2766.6
$ws.Range("K134").Value = 6304.8465
$ws.Range("L134").Value = 8299.799999999999
$ws.Range("M134").Value = -3769.8465
$ws.Range("N134").Value = -13369.8
$ws.Range("H136").Value = 2198.1785
$ws.Range("I136").Value = 1968.25
$ws.Range("J136").Value = 3577.75
$ws.Range("K136").Value = 5904.75
$ws.Range("L136").Value = 10733.25
$ws.Range("M136").Value = -3354.75
$ws.Range("N136").Value = -15833.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 164.8
$ws.Range("I92").Value = 119
$ws.Range("J92").Value = 176.25
$ws.Range("K92").Value = 357
$ws.Range("L92").Value = 528.75
$ws.Range("M92").Value = 891
$ws.Range("N92").Value = -3024.75
$ws.Range("H97").Value = 1313.2142
$ws.Range("I97").Value = 1630.2
$ws.Range("J97").Value = 1137.1111
$ws.Range("K97").Value = 4890.6
$ws.Range("L97").Value = 3411.3333
$ws.Range("M97").Value = -4394.6
$ws.Range("N97").Value = -4403.3333
$ws.Range("H107").Value = 3503133.8
$ws.Range("J107").Value = 5059231.5
$ws.Range("L107").Value = 15177694.5
$ws.Range("N107").Value = -15181534.5
$ws.Range("H113").Value = 1852.579
$ws.Range("I113").Value = 1682
$ws.Range("K113").Value = 5046
$ws.Range("M113").Value = -2876
$ws.Range("H139").Value = 4469.8125
$ws.Range("I139").Value = 1692.9166
$ws.Range("K139").Value = 5078.7498
$ws.Range("M139").Value = 61.2502000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13561.5
$ws.Range("I70").Value = 9946
$ws.Range("J70").Value = 14766.667
$ws.Range("K70").Value = 9946
$ws.Range("L70").Value = 14766.667
$ws.Range("M70").Value = -9676
$ws.Range("N70").Value = -15306.667
$ws.Range("H73").Value = 13561.5
$ws.Range("I73").Value = 9946
$ws.Range("J73").Value = 14766.667
$ws.Range("K73").Value = 9946
$ws.Range("L73").Value = 14766.667
$ws.Range("M73").Value = -9010
$ws.Range("N73").Value = -16638.667
$ws.Range("H80").Value = 3059
$ws.Range("I80").Value = 2422.8572
$ws.Range("J80").Value = 3801.1667
$ws.Range("K80").Value = 2422.8572
$ws.Range("L80").Value = 3801.1667
$ws.Range("M80").Value = -1424.8572
$ws.Range("N80").Value = -5797.1667
$ws.Range("H83").Value = 3059
$ws.Range("I83").Value = 2422.8572
$ws.Range("J83").Value = 3801.1667
$ws.Range("K83").Value = 12114.286
$ws.Range("L83").Value = 19005.8335
$ws.Range("M83").Value = -7122.286
$ws.Range("N83").Value = -28989.8335
$ws.Range("H113").Value = 1160017.2
$ws.Range("I113").Value = 2556.7
$ws.Range("J113").Value = 3089118.2
$ws.Range("K113").Value = 2556.7
$ws.Range("L113").Value = 3089118.2
$ws.Range("M113").Value = -386.6999999999998
$ws.Range("N113").Value = -3093458.2
$ws.Range("H122").Value = 4176.3794
$ws.Range("I122").Value = 3969.9546
$ws.Range("J122").Value = 4825.143
$ws.Range("K122").Value = 11909.8638
$ws.Range("L122").Value = 14475.429
$ws.Range("M122").Value = -9459.863799999999
$ws.Range("N122").Value = -19375.429
$ws.Range("H126").Value = 2305.3125
$ws.Range("I126").Value = 2090.4167
$ws.Range("K126").Value = 6271.250100000001
$ws.Range("M126").Value = -3801.250100000001
$ws.Range("H132").Value = 14687470
$ws.Range("I132").Value = 2201.7
$ws.Range("J132").Value = 63638364
$ws.Range("K132").Value = 6605.099999999999
$ws.Range("L132").Value = 190915092
$ws.Range("M132").Value = -4075.099999999999
$ws.Range("N132").Value = -190920152

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8323.1875
$ws.Range("I7").Value = 7456.1665
$ws.Range("K7").Value = 7456.1665
$ws.Range("M7").Value = -7344.1665
$ws.Range("H16").Value = 1914.6389
$ws.Range("I16").Value = 919.53125
$ws.Range("K16").Value = 919.53125
$ws.Range("M16").Value = -749.53125
$ws.Range("H40").Value = 5799.077
$ws.Range("I40").Value = 5580.727
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 5580.727
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -5444.727
$ws.Range("N40").Value = -7272
$ws.Range("H42").Value = 35416.582
$ws.Range("J42").Value = 32500
$ws.Range("L42").Value = 32500
$ws.Range("N42").Value = -33626
$ws.Range("H49").Value = 35416.582
$ws.Range("J49").Value = 32500
$ws.Range("L49").Value = 32500
$ws.Range("N49").Value = -32794
$ws.Range("H61").Value = 3905.9167
$ws.Range("I61").Value = 3207.375
$ws.Range("J61").Value = 5303
$ws.Range("K61").Value = 3207.375
$ws.Range("L61").Value = 5303
$ws.Range("M61").Value = -3005.375
$ws.Range("N61").Value = -5707
$ws.Range("H113").Value = 3905.9167
$ws.Range("I113").Value = 3207.375
$ws.Range("J113").Value = 5303
$ws.Range("K113").Value = 3207.375
$ws.Range("L113").Value = 5303
$ws.Range("M113").Value = -1037.375
$ws.Range("N113").Value = -9643
$ws.Range("H114").Value = 99999
$ws.Range("J114").Value = 99999
$ws.Range("L114").Value = 99999
$ws.Range("N114").Value = -108677
$ws.Range("H126").Value = 8323.1875
$ws.Range("I126").Value = 7456.1665
$ws.Range("K126").Value = 22368.4995
$ws.Range("M126").Value = -19898.4995
$ws.Range("H132").Value = 3697.524
$ws.Range("J132").Value = 4908.5454
$ws.Range("L132").Value = 14725.6362
$ws.Range("N132").Value = -19785.6362
$ws.Range("H136").Value = 2469.0833
$ws.Range("I136").Value = 2631.3333
$ws.Range("J136").Value = 2306.8333
$ws.Range("K136").Value = 7893.999899999999
$ws.Range("L136").Value = 6920.499899999999
$ws.Range("M136").Value = -5343.999899999999
$ws.Range("N136").Value = -12020.4999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 6502.5
$ws.Range("I6").Value = 6502.5
$ws.Range("K6").Value = 6502.5
$ws.Range("M6").Value = -6387.5
$ws.Range("H8").Value = 1999
$ws.Range("I8").Value = 1999
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1859
$ws.Range("N8").ClearContents()
$ws.Range("H58").Value = 75000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H107").Value = 6863.7856
$ws.Range("I107").Value = 3880.2
$ws.Range("K107").Value = 11640.6
$ws.Range("M107").Value = -9720.599999999999
$ws.Range("H113").Value = 980.7692
$ws.Range("I113").Value = 779.6667
$ws.Range("J113").Value = 1433.25
$ws.Range("K113").Value = 2339.0001
$ws.Range("L113").Value = 4299.75
$ws.Range("M113").Value = -169.0001000000002
$ws.Range("N113").Value = -8639.75
$ws.Range("H126").Value = 3532.25
$ws.Range("I126").Value = 3744.1333
$ws.Range("J126").Value = 2896.6
$ws.Range("K126").Value = 11232.3999
$ws.Range("L126").Value = 8689.799999999999
$ws.Range("M126").Value = -8762.3999
$ws.Range("N126").Value = -13629.8
$ws.Range("H132").Value = 229329.31
$ws.Range("I132").Value = 1882.8611
$ws.Range("K132").Value = 5648.5833
$ws.Range("M132").Value = -3118.5833
$ws.Range("H136").Value = 1431756.4
$ws.Range("I136").Value = 4324.75
$ws.Range("J136").Value = 3334998.8
$ws.Range("K136").Value = 12974.25
$ws.Range("L136").Value = 10004996.4
$ws.Range("M136").Value = -10424.25
$ws.Range("N136").Value = -10010096.4
